$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$partNumber = "71V67603S166PFGI"
$manufacturer = "Renesas"
$supplier = "《此供应商选择了隐藏公司名》"
$colD = "**"
$colE = "**"
$price = "￥19"
$updated = "两年前"
$stock = "<50"
$date = "2023-05-05"

for ($r = 1; $r -le 3; $r++) {
    $ws.Cells.Item($r, 1).Value = $partNumber
    $ws.Cells.Item($r, 2).Value = $manufacturer
    $ws.Cells.Item($r, 3).Value = $supplier
    $ws.Cells.Item($r, 4).Value = $colD
    $ws.Cells.Item($r, 5).Value = $colE
    $ws.Cells.Item($r, 6).Value = $price
    $ws.Cells.Item($r, 7).Value = $updated
    $ws.Cells.Item($r, 8).Value = $stock
    $ws.Cells.Item($r, 10).Value = $false
}

# Column I holds a date-LOOKING string ("2023-05-05") that must stay plain
# text (matches t="inlineStr" in the target), not get auto-coerced into a
# real date serial number the way a normal `.Value = "2023-05-05"` would.
# Stage the text in a scratch range formatted as Text, copy it across, then
# paste-special (values only) into I1:I3 so the destination cells pick up
# the literal string without inheriting the "@" text number format (which
# would otherwise leave a stray style behind).
$scratch = $ws.Range("L1:L3")
$scratch.NumberFormat = "@"
for ($r = 1; $r -le 3; $r++) {
    $ws.Cells.Item($r, 12).Value = $date
}
$scratch.Copy()
$ws.Range("I1:I3").PasteSpecial(-4163)
$scratch.Clear()
